$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value for J10 (new cell with value 5)
$ws.Range("J10").Value = 5

# Set values for G27 and H27 (previously empty)
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 5

# Scroll the frozen pane so the top-left visible cell of the bottom-right pane is C10
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 3
